$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after the "Variant to gene" row (row 12) for the new
# "Gene prioritisation" / "In-vitro experiments" workflow steps. This shifts
# every following row down by two.
$ws.Rows("13:14").Insert()

# --- New row 13: Gene prioritisation ---
$ws.Range("A13").Value = "Gene prioritisation"
$ws.Range("B13").Value = "Gene pathway analysis, protein-protein interaction, gene-variants to bring forward for in-vitro experiments"

# --- New row 14: In-vitro experiments ---
$ws.Range("A14").Value = "In-vitro experiments"

# Notes column for the new "Gene prioritisation" row.
$ws.Range("E13").Value = "Create hypothesis of variant-gene mechanism for in-vitro experiments"

$ws.Range("B14").Value = "CRISPR-Cas9 experiments"

# The inserted row 14 picked up a stray (empty) Notes cell from row 12's
# layout; row 14 shouldn't have a Notes column entry at all.
$ws.Range("E14").Clear()

# Copy the "TO BE DONE" formatting (from the row that used to sit at 14, now
# at 16) onto the two freshly inserted rows' status columns.
$ws.Range("C16:D16").Copy()
$ws.Range("C13:D13").PasteSpecial(-4122)
$ws.Range("C14:D14").PasteSpecial(-4122)
$ws.Range("C13:D13").Value = "TO BE DONE"
$ws.Range("C14:D14").Value = "TO BE DONE"

# Row 13 wraps onto multiple lines once filled in.
$ws.Rows(13).RowHeight = 43.5

# --- Row 12: the "Variant to gene" analysis is now finished and reported ---
$ws.Range("B12").Value = "Variant to gene mapping with 8 line of evidence"
$ws.Range("C2:D2").Copy()
$ws.Range("C12:D12").PasteSpecial(-4122)
$ws.Range("C12:D12").Value = "DONE"
$ws.Range("E12").Clear()

# Row 15 (the old row 13, shifted down) keeps an empty Notes cell instead of
# the stale "Need to draft an analysis plan" note.
$ws.Range("E15").ClearContents()

$ws.Range("B12").Select()
